# Change the "type" field for the M_FOL_date and M_time rows on the
# "survey" sheet from "date" / "time" to "text", per the commit message
# ("Changed date and time fields for JGI app and verified database
# persistence"). This removes the now-unused "date" and "time" shared
# strings automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Row 2 -> field M_FOL_date: type date -> text
$ws.Range("C2").Value = "text"

# Row 4 -> field M_time: type time -> text
$ws.Range("C4").Value = "text"

# Update the active selection on the survey sheet to C6
$ws.Range("C6").Select()
